$d = $word.ActiveDocument

# 1) Merge the split "Date:" runs into a single run's text.
#    The cell currently contains separate runs: "Date:", "06", " 0", "5", " 2025"
#    which render as "Date:06 05 2025". Collapse them into one run with that text
#    by deleting the runs after the first and setting the first run's text.
$table = $d.Tables.Item(1)
$dateCell = $table.Cell(1, 2)
$dateRange = $dateCell.Range
$dateRange.End = $dateRange.End - 1
$dateRange.Text = "Date:06 05 2025"

# 2) Fill in the empty run in the "Evaluation of performance" answer cell
#    (row 3, column 2) with the new reflection text.
$evalCell = $table.Cell(3, 2)
$evalRange = $evalCell.Range
$evalRange.End = $evalRange.End - 1
$evalRange.Text = "I have prepared for the presentation by adding screenshots of development build and sharing them to the team to add to posters and powerpoints. As well as this I assisted in the making of the powerpoint in a long meeting with jeremy and heytham. In this we created the presentation structure and content. We had a further one where we discussed who had responsibility of presenting which part, I had it over the live demo. We made notes and discussed our parts further. We then done the same on the day of the presentation. Th presentation itself went really well I thought, I put my all into to conveying the ideas we had developed over the time of this project as well as making it engaging for the clients. I am immensely happy with the effort but In by the 4 of us who were present that day. I involved the clients in the demo as well as trying to show the good work we had all put into developing the soloution."
$evalRange.Font.Name = $evalRange.Font.Name
